# Basic Prototype Ready! (Team Commit #2)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update existing Week 9 breakdown cells (B10:B12) ---
$ws.Range("B10").Value = "TimeSheet UI"
$ws.Range("B11").Value = "Permissions"
$ws.Range("B12").Value = "Data Export"

# --- Add new small assignment table at rows 79-83 ---
$ws.Range("B79").Value = "Item"
$ws.Range("C79").Value = "Arnav"
$ws.Range("D79").Value = "Sid"
$ws.Range("E79").Value = "Lohit"

$ws.Range("A80").Value = "Wed"
$ws.Range("B80").Value = "Basic Setup"
$ws.Range("C80").Value = "Frontend "
$ws.Range("D80").Value = "Django"
$ws.Range("E80").Value = "Django"

$ws.Range("A81").Value = "Thu"
$ws.Range("B81").Value = "Basic Setup"
$ws.Range("C81").Value = "Frontend "
$ws.Range("D81").Value = "Django"
$ws.Range("E81").Value = "Django"

$ws.Range("A82").Value = "Fri"
$ws.Range("B82").Value = "Frontend"
$ws.Range("C82").Value = "Data Collection"
$ws.Range("D82").Value = "Updation in Activity Report"
$ws.Range("E82").Value = "Choices API Integration"

$ws.Range("A83").Value = "Sat"
$ws.Range("B83").Value = "Frontend"
$ws.Range("C83").Value = "Data Collection"
$ws.Range("D83").Value = "Updation in Activity Report"
$ws.Range("E83").Value = "Choices API Integration"

# --- Update view state to reflect scroll position / selection ---
try {
    $excel.ActiveWindow.ScrollRow = 73
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll-position is cosmetic view state only; ignore if unsupported.
}
$ws.Range("A85").Select()
